# player_stats.xlsx edit:
# Append 19 new player-stat rows (48-66) to the sheet. These duplicate
# the existing "KAGS #7158" stat block already present in rows 4-47,
# except that row 52's player tag is missing the space
# ("KAGS#7158" instead of "KAGS #7158"), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt 19; $i++) {
    $destRow = 48 + $i

    $ws.Cells.Item($destRow, 1).Value  = "KAGS #7158"
    $ws.Cells.Item($destRow, 2).Value  = 1.01
    $ws.Cells.Item($destRow, 3).Value  = 57
    $ws.Cells.Item($destRow, 4).Value  = "Jett"
    $ws.Cells.Item($destRow, 5).Value  = 23.3
    $ws.Cells.Item($destRow, 6).Value  = 14
    $ws.Cells.Item($destRow, 7).Value  = 138
    $ws.Cells.Item($destRow, 8).Value  = 156
    $ws.Cells.Item($destRow, 9).Value  = 0
    $ws.Cells.Item($destRow, 10).Value = 15
    $ws.Cells.Item($destRow, 11).Value = "Nickel"
    $ws.Cells.Item($destRow, 12).Value = "['Rusher', 'Straight Up Winner']"
}

# Row 52 has the player tag without the space between "KAGS" and "#7158".
$ws.Cells.Item(52, 1).Value = "KAGS#7158"
